$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember the width of the column that will become the left neighbour of
# the freshly inserted column so the new column can inherit it, mirroring
# Excel's native "Insert Column" behaviour.
$existingWidth = $ws.Columns.Item(13).ColumnWidth

# Insert a new, blank column before column N (14) - this shifts the old
# N,O,P data right into O,P,Q.
$null = $ws.Columns.Item(14).Insert()
$ws.Columns.Item(14).ColumnWidth = $existingWidth

# Make "Repayment schedule" the active sheet and select cell K13 on it.
$ws.Activate()
$null = $ws.Range("K13").Select()
